$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 2..151) holds image filenames like "counting/counting_0.png".
# Strip the "counting/" directory prefix so they become "counting_0.png".
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "counting/counting_*.png") {
        $cell.Value2 = $val -replace "^counting/", ""
    }
}
